$wb = $excel.ActiveWorkbook

# --- Add the new "SearchSkill" sheet right before "Register" ---
$registerSheet = $wb.Worksheets.Item("Register")
$searchSkill = $wb.Worksheets.Add($registerSheet)
$searchSkill.Name = "SearchSkill"

$searchSkill.Range("A1").Value = "SearchCategory"
$searchSkill.Range("A2").Value = "Fun & Lifestyle"
$searchSkill.Range("A1:A2").NumberFormat = "General"

$searchSkill.Columns.Item(1).ColumnWidth = 13.2857142857

$searchSkill.Range("C2").Select()

# --- ShareSkill sheet view + data tweaks ---
$shareSkill = $wb.Worksheets.Item("ShareSkill")
$shareSkill.Range("M8").Select()

$shareSkill.Range("L2").Value = 44910
$shareSkill.Range("M2").Value = 45275

$shareSkill.Range("L3").Value = 44917
$shareSkill.Range("M3").Value = 45282

$shareSkill.Range("L4").Value = 44920
$shareSkill.Range("M4").Value = 45285

$shareSkill.Range("L5").Value = 44793
$shareSkill.Range("M5").Value = 45158

# Re-select SearchSkill so it ends up the active tab, matching the target.
$searchSkill.Select()
$searchSkill.Range("C2").Select()
